# Add a third worksheet "ODI Batting Extra" after the existing "ODI Batting"
# sheet, carrying per-match batting extras (boundary counts, % of team runs,
# man-of-the-match flag), mirroring the existing two sheets' layout/style.

$wb = $excel.ActiveWorkbook

# --- create + place the new sheet at the end of the tab strip -------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# --- header row -------------------------------------------------------------
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Match the bold / bordered / centered header styling already used on the
# "Player Info" and "ODI Batting" sheets by copying the format of an
# existing header cell onto the whole new header row (keeps the same
# underlying style, rather than inventing a new one).
$headerFormat = $wb.Worksheets.Item(1).Range("A1")
$headerFormat.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- data rows ---------------------------------------------------------------
# Text-typed columns (MATCH_CODE, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL,
# MAN_OF_MATCH) hold values that look numeric ("4472", "7", "0", "40.85%"),
# so force text formatting before assigning and reset the style back to the
# plain/default one afterwards (matches the unstyled data cells elsewhere).
$textCells = @(
    @{ Cell = "A2"; Value = "4472" },
    @{ Cell = "C2"; Value = "7" },
    @{ Cell = "D2"; Value = "0" },
    @{ Cell = "E2"; Value = "40.85%" },
    @{ Cell = "F2"; Value = "NO" },

    @{ Cell = "A3"; Value = "4473" },
    @{ Cell = "C3"; Value = "0" },
    @{ Cell = "D3"; Value = "0" },
    @{ Cell = "E3"; Value = "" },
    @{ Cell = "F3"; Value = "NO" },

    @{ Cell = "A4"; Value = "4476" },
    @{ Cell = "C4"; Value = "7" },
    @{ Cell = "D4"; Value = "0" },
    @{ Cell = "E4"; Value = "11.75%" },
    @{ Cell = "F4"; Value = "NO" }
)

foreach ($entry in $textCells) {
    $rng = $ws.Range($entry.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $entry.Value
    $rng.Style = "Normal"
}

# BATTING_POSITION is a genuine number (3) for every row.
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 3
